$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three fixed answer strings to cycle through for columns B (room letter) and D (answer text)
$letters = @("א", "ב", "ג")
$answers = @(
    "דני צודק,צודק,צדק,כן,נכון",
    "מיכל צודקת,כן,צודק,נכון",
    "דני לא צודק,דני טועה,טעות,לא,טואה,תואה,תאות,לא צודק,לא נכון,דני שוגה"
)

$row = 114
for ($room = 792; $room -le 802; $room++) {
    for ($i = 0; $i -lt 3; $i++) {
        $ws.Cells.Item($row, 1).Value = $room
        $ws.Cells.Item($row, 2).Value = $letters[$i]
        $ws.Cells.Item($row, 4).Value = $answers[$i]
        $row++
    }
}

$ws.Range("D143").Select()
